$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    86.36363636363636,
    87.991858887381284,
    86.160108548168239,
    88.331071913161466,
    88.331071913161466,
    88.263229308005435,
    87.516960651289011,
    85.549525101763919,
    86.567164179104466,
    86.635006784260511,
    86.092265943012208,
    86.499321573948436,
    88.12754409769336,
    86.770691994572587,
    88.12754409769336,
    87.516960651289011,
    88.19538670284939,
    87.449118046132966,
    88.19538670284939,
    88.19538670284939,
    88.059701492537314,
    87.788331071913163,
    87.720488466757118,
    87.652645861601087,
    87.109905020352784,
    88.263229308005435,
    88.331071913161466,
    87.924016282225239,
    87.856173677069208,
    88.059701492537314,
    88.19538670284939,
    88.12754409769336,
    88.331071913161466,
    87.584803256445042,
    87.516960651289011,
    87.24559023066486,
    87.584803256445042,
    87.652645861601087,
    87.788331071913163,
    87.991858887381284,
    88.263229308005435,
    88.19538670284939,
    88.059701492537314,
    87.924016282225239,
    87.991858887381284,
    87.516960651289011,
    86.024423337856177,
    88.19538670284939
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Re-assign A1 so a fresh (new) shared string entry gets created for it,
# matching the new-model/new-dataset string table growth.
$ws.Range("A1").Value = "HK_R_acc_LG"
